$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds free-form text (values like "68.123.45" use dots as
# thousands separators, not decimals), so force text formatting before writing
# the new values to avoid Excel re-interpreting them as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.203.00"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "2.639.57"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "597.13"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "155.17"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "2.638.21"
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("D10").Value = "0.145"
$ws.Range("E10").Value = "  +8.37%  "
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("D14").Value = "27.99"
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("D15").Value = "0.0000192"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("D16").Value = "3.119.78"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "68.069.41"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "2.648.89"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "11.37"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").Value = "363.25"
$ws.Range("E20").Value = "  -1.15%  "
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").Value = "4.37"
$ws.Range("E22").Value = "  +3.41%  "
$ws.Range("D23").Value = "4.83"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").Value = "2.06"
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("D25").Value = "75.02"
$ws.Range("E25").Value = "  +3.28%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("D27").Value = "9.69"
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("E28").Value = "  +1.83%  "
$ws.Range("D29").Value = "2.774.52"
$ws.Range("E29").Value = "  +0.87%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").Value = "557.67"
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("D32").Value = "8.00"
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("E35").Value = "  +1.55%  "
$ws.Range("E36").Value = "  -0.20%  "
$ws.Range("E37").Value = "  +3.36%  "
$ws.Range("D38").Value = "161.16"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").Value = "19.30"
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("D40").Value = "0.372"
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("D41").Value = "1.88"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").Value = "5.31"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("E43").Value = "  +4.63%  "
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("D45").Value = "2.63"
$ws.Range("E45").Value = "  -0.37%  "
$ws.Range("D47").Value = "40.44"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("D48").Value = "158.69"
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("D49").Value = "3.74"
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("D50").Value = "21.94"
$ws.Range("E50").Value = "  +0.82%  "
$ws.Range("D51").Value = "0.0784"
$ws.Range("E51").Value = "  +0.95%  "
